# Fruta / hortaliza, semanal
# Two new weekly price rows are inserted into the "Limón" daily-logic
# subconjunto sheet. They land right after the existing 2021-08-30
# ("1a amarillo" / "2a amarillo") rows, which pushes every following
# row down by two positions (row 515 -> 517, ... row 595 -> 597) while
# keeping their original data intact. The dimension grows from
# A1:T595 to A1:T597.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 515:516 - everything from old row 515 downward
# shifts down to 517 onward, carrying its original values with it.
$ws.Rows("515:516").Insert()

# Seed the two new rows with the same static/lookup columns (A-L, Q, R, T)
# as the rows that used to occupy 515/516 (now at 517/518), since those
# columns are unchanged across the whole "Femacal de La Calera" / "Limón"
# block. The numeric columns (D, M, N, O, P, S) are overwritten right after
# with the real new weekly figures.
$ws.Range("A517:T517").Copy() | Out-Null
$ws.Range("A515:T515").PasteSpecial() | Out-Null

$ws.Range("A518:T518").Copy() | Out-Null
$ws.Range("A516:T516").PasteSpecial() | Out-Null

$excel.CutCopyMode = 0

# Row 515: 1a amarillo, week of 2021-10-05 (serial 44474)
$ws.Range("D515").Value = 44474
$ws.Range("M515").Value = 295
$ws.Range("N515").Value = 4000
$ws.Range("O515").Value = 4500
$ws.Range("P515").Value = 4236
$ws.Range("S515").Value = 265

# Row 516: 2a amarillo, week of 2021-10-05 (serial 44474)
$ws.Range("D516").Value = 44474
$ws.Range("M516").Value = 177
$ws.Range("N516").Value = 3000
$ws.Range("O516").Value = 3500
$ws.Range("P516").Value = 3246
$ws.Range("S516").Value = 203
